# "create project - ultimos casos de prueba"
# Adds two new worksheets (Hoja15, Hoja16) with "Rol" test data after Hoja14,
# tweaks Hoja14's selection (no longer the active/selected tab), and leaves
# Hoja16 as the newly active/selected sheet.

$wb = $excel.ActiveWorkbook

# Header / row values shared by both new sheets.
$headers = @("Username", "Contraseña", "nombre proyecto", "Descripcion", "Fecha de inicio", "Rol", "member")
$rowCommon = @("pepeusername", "P4ssword.", "desctest1", "Length is 12", "13 October 2021", "Team Member")

# ---- Hoja15 -----------------------------------------------------------
$hoja14 = $wb.Worksheets.Item("Hoja14")
$hoja15 = $wb.Worksheets.Add($null, $hoja14)
$hoja15.Name = "Hoja15"

$hoja15.Range("A1").Value = $headers[0]
$hoja15.Range("B1").Value = $headers[1]
$hoja15.Range("C1").Value = $headers[2]
$hoja15.Range("D1").Value = $headers[3]
$hoja15.Range("E1").Value = $headers[4]
$hoja15.Range("F1").Value = $headers[5]
$hoja15.Range("G1").Value = $headers[6]

$hoja15.Range("A2").Value = $rowCommon[0]
$hoja15.Range("B2").Value = $rowCommon[1]
$hoja15.Range("C2").Value = $rowCommon[2]
$hoja15.Range("D2").Value = $rowCommon[3]
$hoja15.Range("E2").Value = $rowCommon[4]
$hoja15.Range("F2").Value = $rowCommon[5]
$hoja15.Range("G2").Value = "asdf123"

$hoja15.Range("A3").Value = $rowCommon[0]
$hoja15.Range("B3").Value = $rowCommon[1]
$hoja15.Range("C3").Value = $rowCommon[2]
$hoja15.Range("D3").Value = $rowCommon[3]
$hoja15.Range("E3").Value = $rowCommon[4]
$hoja15.Range("F3").Value = $rowCommon[5]
$hoja15.Range("G3").Value = "qwer123"

$hoja15.Range("A1:G3").Select() | Out-Null

# ---- Hoja16 -----------------------------------------------------------
$hoja16 = $wb.Worksheets.Add($null, $hoja15)
$hoja16.Name = "Hoja16"

$hoja16.Range("A1").Value = $headers[0]
$hoja16.Range("B1").Value = $headers[1]
$hoja16.Range("C1").Value = $headers[2]
$hoja16.Range("D1").Value = $headers[3]
$hoja16.Range("E1").Value = $headers[4]
$hoja16.Range("F1").Value = $headers[5]
$hoja16.Range("G1").Value = $headers[6]

$hoja16.Range("A2").Value = $rowCommon[0]
$hoja16.Range("B2").Value = $rowCommon[1]
$hoja16.Range("C2").Value = $rowCommon[2]
$hoja16.Range("D2").Value = $rowCommon[3]
$hoja16.Range("E2").Value = $rowCommon[4]
$hoja16.Range("F2").Value = $rowCommon[5]
$hoja16.Range("G2").Value = "jositom"

$hoja16.Range("A3").Value = $rowCommon[0]
$hoja16.Range("B3").Value = $rowCommon[1]
$hoja16.Range("C3").Value = $rowCommon[2]
$hoja16.Range("D3").Value = $rowCommon[3]
$hoja16.Range("E3").Value = $rowCommon[4]
$hoja16.Range("F3").Value = $rowCommon[5]
$hoja16.Range("G3").Value = "jonam"

# Hoja14 keeps its data but is no longer the selected/active tab - select
# a fresh range on it to mirror the recorded final state.
$hoja14.Range("A1:E2").Select() | Out-Null

# Hoja16 ends up the active sheet/tab with G3 selected.
$hoja16.Activate() | Out-Null
$hoja16.Range("G3").Select() | Out-Null
